# --- Sheet "Input": rebuild with the new standard template column layout ---
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Wipe the old layout (old headers, old style, old columns) completely.
$ws.Cells.Clear()

$headers = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$rows = @(
    @("2025-09-15","2025-10-08","이노메탈","이노메탈@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","B. 도장","2월 청구","KS규격-1",930,1900,1943700,"제이비엔지니어링"),
    @("2025-09-13","2025-09-26","이노메탈","이노메탈@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","5. 운반비","일반자재","기타","운반비","KS규격-2",1,0,0,$null),
    @("2025-09-01","2025-08-30","이노메탈","이노메탈@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","B. 도장","2월 청구","KS규격-3",2010,1500,3316500,"제이비엔지니어링"),
    @("2025-08-26","2025-09-16","이노메탈","이노메탈@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","5. 운반비","일반자재","기타","4월 운반비","KS규격-4",1,0,0,$null),
    @("2025-08-24","2025-10-09","이노메탈","이노메탈@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","B. 도장","도장3차","KS규격-5",520,1500,858000,$null),
    @("2025-09-14","2025-10-09","이노메탈","이노메탈@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","B. 도장","품목명 없음","KS규격-6",1307,2600,3738020,$null)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2

    # Columns A (발주일자) and B (납기일자) hold plain date-strings, not real
    # Excel dates - force text format per-cell first so assigning a
    # "yyyy-mm-dd" string doesn't get silently reinterpreted as a date serial.
    $ws.Cells.Item($excelRow, 1).NumberFormat = "@"
    $ws.Cells.Item($excelRow, 2).NumberFormat = "@"

    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($excelRow, $c + 1).Value = $val
        }
    }
}

Write-Host "Input sheet rebuilt"
